$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 1364
$ws.Range("E2").Value = 205
$ws.Range("F2").Value = 205
$ws.Range("G2").Value = 120
$ws.Range("H2").Value = 109
$ws.Range("I2").Value = 108
$ws.Range("J2").Value = 2
$ws.Range("K2").Value = 2106
$ws.Range("L2").Value = 929
$ws.Range("M2").Value = 1177
$ws.Range("N2").Value = 1171
$ws.Range("O2").Value = 6
$ws.Range("P2").Value = 148
$ws.Range("Q2").Value = 376
$ws.Range("R2").Value = -105
$ws.Range("S2").Value = -237
$ws.Range("T2").Value = 63
$ws.Range("U2").Value = 314
$ws.Range("V2").Value = 507
$ws.Range("W2").Value = 15.05
$ws.Range("X2").Value = 8.02
$ws.Range("Y2").Value = 10.4
$ws.Range("Z2").Value = 5.42
$ws.Range("AA2").Value = 78.95
$ws.Range("AB2").Value = 691.48
$ws.Range("AC2").Value = 363
$ws.Range("AD2").Value = 25.6
$ws.Range("AE2").Value = 3947
$ws.Range("AF2").Value = 2.36
$ws.Range("AG2").Value = 200
$ws.Range("AH2").Value = 2.15
$ws.Range("AI2").Value = 55.06
$ws.Range("AJ2").Value = 29672700
$ws.Range("D3").Value = 1577
$ws.Range("E3").Value = 290
$ws.Range("F3").Value = 290
$ws.Range("G3").Value = 242
$ws.Range("H3").Value = 217
$ws.Range("I3").Value = 214
$ws.Range("J3").Value = 3
$ws.Range("K3").Value = 2135
$ws.Range("L3").Value = 808
$ws.Range("M3").Value = 1328
$ws.Range("N3").Value = 1319
$ws.Range("O3").Value = 9
$ws.Range("P3").Value = 148
$ws.Range("Q3").Value = 465
$ws.Range("R3").Value = -178
$ws.Range("S3").Value = -215
$ws.Range("T3").Value = 99
$ws.Range("U3").Value = 366
$ws.Range("V3").Value = 353
$ws.Range("W3").Value = 18.38
$ws.Range("X3").Value = 13.74
$ws.Range("Y3").Value = 17.17
$ws.Range("Z3").Value = 10.22
$ws.Range("AA3").Value = 60.82
$ws.Range("AB3").Value = 791.14
$ws.Range("AC3").Value = 720
$ws.Range("AD3").Value = 28.46
$ws.Range("AE3").Value = 4445
$ws.Range("AF3").Value = 4.61
$ws.Range("AG3").Value = 220
$ws.Range("AH3").Value = 1.07
$ws.Range("AI3").Value = 30.54
$ws.Range("AJ3").Value = 29672700
$ws.Range("D4").Value = 1768
$ws.Range("E4").Value = 384
$ws.Range("F4").Value = 384
$ws.Range("G4").Value = 356
$ws.Range("H4").Value = 282
$ws.Range("I4").Value = 282
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 2365
$ws.Range("L4").Value = 832
$ws.Range("M4").Value = 1533
$ws.Range("N4").Value = 1524
$ws.Range("O4").Value = 9
$ws.Range("P4").Value = 148
$ws.Range("Q4").Value = 496
$ws.Range("R4").Value = -241
$ws.Range("S4").Value = -101
$ws.Range("T4").Value = 50
$ws.Range("U4").Value = 446
$ws.Range("V4").Value = 317
$ws.Range("W4").Value = 21.73
$ws.Range("X4").Value = 15.97
$ws.Range("Y4").Value = 19.82
$ws.Range("Z4").Value = 12.55
$ws.Range("AA4").Value = 54.27
$ws.Range("AB4").Value = 929.3200000000001
$ws.Range("AC4").Value = 949
$ws.Range("AD4").Value = 22.6
$ws.Range("AE4").Value = 5135
$ws.Range("AF4").Value = 4.18
$ws.Range("AG4").Value = 290
$ws.Range("AH4").Value = 1.35
$ws.Range("AI4").Value = 30.55
$ws.Range("AJ4").Value = 29672700
$ws.Range("D5").Value = 2044
$ws.Range("E5").Value = 505
$ws.Range("F5").Value = 505
$ws.Range("G5").Value = 482
$ws.Range("H5").Value = 397
$ws.Range("I5").Value = 392
$ws.Range("J5").Value = 5
$ws.Range("K5").Value = 2743
$ws.Range("L5").Value = 919
$ws.Range("M5").Value = 1824
$ws.Range("N5").Value = 1810
$ws.Range("O5").Value = 15
$ws.Range("P5").Value = 148
$ws.Range("Q5").Value = 501
$ws.Range("R5").Value = -236
$ws.Range("S5").Value = -115
$ws.Range("T5").Value = 332
$ws.Range("U5").Value = 169
$ws.Range("V5").Value = 289
$ws.Range("W5").Value = 24.71
$ws.Range("X5").Value = 19.41
$ws.Range("Y5").Value = 23.49
$ws.Range("Z5").Value = 15.53
$ws.Range("AA5").Value = 50.36
$ws.Range("AB5").Value = 1122.51
$ws.Range("AC5").Value = 1319
$ws.Range("AD5").Value = 25.24
$ws.Range("AE5").Value = 6098
$ws.Range("AF5").Value = 5.46
$ws.Range("AG5").Value = 400
$ws.Range("AH5").Value = 1.2
$ws.Range("AI5").Value = 30.32
$ws.Range("AJ5").Value = 29672700
$ws.Range("D6").Value = 2269
$ws.Range("E6").Value = 540
$ws.Range("F6").Value = 540
$ws.Range("G6").Value = 561
$ws.Range("H6").Value = 425
$ws.Range("I6").Value = 421
$ws.Range("K6").Value = 2787
$ws.Range("L6").Value = 965
$ws.Range("M6").Value = 1822
$ws.Range("N6").Value = 1794
$ws.Range("P6").Value = 148
$ws.Range("Q6").Value = 621
$ws.Range("R6").Value = -200
$ws.Range("S6").Value = -506
$ws.Range("T6").Value = 183
$ws.Range("U6").Value = 439
$ws.Range("V6").Value = 186
$ws.Range("W6").Value = 23.8
$ws.Range("X6").Value = 18.73
$ws.Range("Y6").Value = 23.36
$ws.Range("Z6").Value = 15.38
$ws.Range("AA6").Value = 53
$ws.Range("AB6").Value = 1328.77
$ws.Range("AC6").Value = 1418
$ws.Range("AD6").Value = 36.59
$ws.Range("AE6").Value = 6175
$ws.Range("AF6").Value = 8.4
$ws.Range("AG6").Value = 450
$ws.Range("AH6").Value = 0.87
$ws.Range("AI6").Value = 31.06
$ws.Range("AJ6").Value = 29672700
$ws.Range("D7").Value = 2586
$ws.Range("E7").Value = 644
$ws.Range("G7").Value = 674
$ws.Range("H7").Value = 510
$ws.Range("I7").Value = 508
$ws.Range("K7").Value = 4536
$ws.Range("L7").Value = 1870
$ws.Range("M7").Value = 2666
$ws.Range("N7").Value = 2438
$ws.Range("P7").Value = 151
$ws.Range("Q7").Value = 625
$ws.Range("R7").Value = -1626
$ws.Range("S7").Value = 1157
$ws.Range("T7").Value = 1073
$ws.Range("U7").Value = -594
$ws.Range("W7").Value = 24.9
$ws.Range("X7").Value = 19.74
$ws.Range("Y7").Value = 24.02
$ws.Range("Z7").Value = 13.94
$ws.Range("AA7").Value = 70.16
$ws.Range("AC7").Value = 1675
$ws.Range("AD7").Value = 54.26
$ws.Range("AE7").Value = 7900
$ws.Range("AF7").Value = 11.51
$ws.Range("AG7").Value = 485
$ws.Range("AH7").Value = 0.53
$ws.Range("AI7").Value = 28.32
$ws.Range("D8").Value = 3078
$ws.Range("E8").Value = 806
$ws.Range("G8").Value = 805
$ws.Range("H8").Value = 629
$ws.Range("I8").Value = 624
$ws.Range("K8").Value = 7147
$ws.Range("L8").Value = 3250
$ws.Range("M8").Value = 3897
$ws.Range("N8").Value = 3813
$ws.Range("P8").Value = 157
$ws.Range("Q8").Value = 920
$ws.Range("R8").Value = -368
$ws.Range("S8").Value = -99
$ws.Range("T8").Value = 130
$ws.Range("U8").Value = 750
$ws.Range("W8").Value = 26.18
$ws.Range("X8").Value = 20.43
$ws.Range("Y8").Value = 18.6
$ws.Range("Z8").Value = 9.65
$ws.Range("AA8").Value = 83.39
$ws.Range("AC8").Value = 1983
$ws.Range("AD8").Value = 45.49
$ws.Range("AE8").Value = 12358
$ws.Range("AF8").Value = 7.3
$ws.Range("AG8").Value = 494
$ws.Range("AH8").Value = 0.55
$ws.Range("AI8").Value = 23.5
$ws.Range("D9").Value = 3654
$ws.Range("E9").Value = 996
$ws.Range("G9").Value = 990
$ws.Range("H9").Value = 775
$ws.Range("I9").Value = 766
$ws.Range("K9").Value = 7710
$ws.Range("L9").Value = 3229
$ws.Range("M9").Value = 4481
$ws.Range("N9").Value = 4438
$ws.Range("P9").Value = 156
$ws.Range("Q9").Value = 971
$ws.Range("R9").Value = -406
$ws.Range("S9").Value = -94
$ws.Range("T9").Value = 112
$ws.Range("U9").Value = 822
$ws.Range("W9").Value = 27.27
$ws.Range("X9").Value = 21.22
$ws.Range("Y9").Value = 18.57
$ws.Range("Z9").Value = 10.44
$ws.Range("AA9").Value = 72.06999999999999
$ws.Range("AC9").Value = 2434
$ws.Range("AD9").Value = 37.06
$ws.Range("AE9").Value = 14382
$ws.Range("AF9").Value = 6.27
$ws.Range("AG9").Value = 530
$ws.Range("AH9").Value = 0.59
$ws.Range("AI9").Value = 20.52
